$d = $word.ActiveDocument

# --- 1. Append " - KTPM1" (en dash) after the Group ID number run ---
$dash = [string][char]0x2013
$suffix = " " + $dash + " KTPM1"

$findRng = $d.Content
$null = $findRng.Find.Execute("Group ID: 6", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findRng.Collapse(0)   # wdCollapseEnd
$findRng.InsertAfter($suffix)

# Grab a fresh range over the inserted text and apply the same look as the
# "6" run (Arial / bold / blue 0070C0 / szCs 18)
$newRng = $d.Content
$null = $newRng.Find.Execute($dash + " KTPM1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newRng.MoveStart(1, -1)  # wdCharacter - extend one char left to capture the leading space
$newRng.Font.NameAscii = "Arial"
$newRng.Font.NameOther = "Arial"
$newRng.Font.NameBi = "Arial"
$newRng.Font.Bold = $true
$newRng.Font.Color = 12611584
$newRng.Font.SizeBi = 9
